# Hortaliza, Vega Central Mapocho de Santiago - Alcachofa
# Weekly update: insert two new report rows (a new "Española" price
# observation, qualities Primera/Segunda) at the top of the data block,
# pushing the existing rows (261-286) down to (263-288).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 261; everything below
# (the existing rows 261:286) shifts down to 263:288.
$ws.Rows("261:262").Insert()

# --- New row 261 -----------------------------------------------------
$ws.Range("A261").Value = 9
$ws.Range("B261").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C261").Value = "Metropolitana"
$ws.Range("D261").Value = 44461
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 100112013
$ws.Range("G261").Value = "Alcachofa"
$ws.Range("H261").Value = "Española"
$ws.Range("I261").Value = "Primera"
$ws.Range("J261").Value = 52
$ws.Range("K261").Value = 12000
$ws.Range("L261").Value = 13000
$ws.Range("M261").Value = 12500
$ws.Range("N261").Value = "`$/caja 30 unidades"
$ws.Range("O261").Value = "Provincia de Limarí"
$ws.Range("P261").Value = 417
$ws.Range("Q261").Value = 30
$ws.Range("R261").Value = "Hortaliza"

# --- New row 262 -----------------------------------------------------
$ws.Range("A262").Value = 9
$ws.Range("B262").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C262").Value = "Metropolitana"
$ws.Range("D262").Value = 44461
$ws.Range("E262").Value = 13
$ws.Range("F262").Value = 100112013
$ws.Range("G262").Value = "Alcachofa"
$ws.Range("H262").Value = "Española"
$ws.Range("I262").Value = "Segunda"
$ws.Range("J262").Value = 25
$ws.Range("K262").Value = 10000
$ws.Range("L262").Value = 11000
$ws.Range("M262").Value = 10480
$ws.Range("N262").Value = "`$/caja 40 unidades"
$ws.Range("O262").Value = "Provincia de Limarí"
$ws.Range("P262").Value = 262
$ws.Range("Q262").Value = 40
$ws.Range("R262").Value = "Hortaliza"

# Make sure the date cells keep the workbook's date/time display format
# (same custom numFmt already used by every other row in column D).
$ws.Range("D261").NumberFormat = $ws.Range("D263").NumberFormat
$ws.Range("D262").NumberFormat = $ws.Range("D263").NumberFormat
